# Refresh the crypto snapshot: latest Price (col D) and Volume(1h) % (col E)
# for each coin row, plus a roster change in row 51 (BabyDogeCoin -> USDD).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.931.81'
$ws.Range('E2').Value = '  +0.69%  '

$ws.Range('D3').Value = '1.551.55'
$ws.Range('E3').Value = '  +0.89%  '

$ws.Range('E4').Value = '  +0.47%  '

$ws.Range('D5').Formula = "'" + '206.82'
$ws.Range('E5').Value = '  +0.61%  '

$ws.Range('E6').Value = '  +1.13%  '

$ws.Range('E7').Value = '  +0.35%  '

$ws.Range('E8').Value = '  +1.73%  '

$ws.Range('D9').Formula = "'" + '21.67'
$ws.Range('E9').Value = '  +2.47%  '

$ws.Range('E10').Value = '  +1.66%  '

$ws.Range('D11').Formula = "'" + '0.0858'
$ws.Range('E11').Value = '  +0.57%  '

$ws.Range('D12').Value = '1.772.42'
$ws.Range('E12').Value = '  +0.51%  '

$ws.Range('D13').Value = '1.548.69'
$ws.Range('E13').Value = '  +0.70%  '

$ws.Range('D14').Formula = "'" + '3.72'
$ws.Range('E14').Value = '  +1.49%  '

$ws.Range('E15').Value = '  +1.95%  '

$ws.Range('D16').Formula = "'" + '61.79'
$ws.Range('E16').Value = '  +1.33%  '

$ws.Range('D17').Value = '26.917.55'
$ws.Range('E17').Value = '  +0.83%  '

$ws.Range('D18').Formula = "'" + '215.58'
$ws.Range('E18').Value = '  +1.90%  '

$ws.Range('D19').Value = '0.0₃0688'
$ws.Range('E19').Value = '  +0.26%  '

$ws.Range('E20').Value = '  +0.16%  '

$ws.Range('E21').Value = '  +0.48%  '

$ws.Range('D22').Formula = "'" + '4.04'
$ws.Range('E22').Value = '  +0.09%  '

$ws.Range('D23').Formula = "'" + '9.14'
$ws.Range('E23').Value = '  +1.64%  '

$ws.Range('D24').Formula = "'" + '1.99'
$ws.Range('E24').Value = '  +0.17%  '

$ws.Range('D25').Formula = "'" + '152.62'
$ws.Range('E25').Value = '  -0.04%  '

$ws.Range('E26').Value = '  +3.32%  '

$ws.Range('E27').Value = '  +0.32%  '

$ws.Range('E28').Value = '  +0.40%  '

$ws.Range('E29').Value = '  +1.02%  '

$ws.Range('E30').Value = '  +1.93%  '

$ws.Range('E31').Value = '  -0.44%  '

$ws.Range('E32').Value = '  +0.00%  '

$ws.Range('D33').Value = '1.416.14'
$ws.Range('E33').Value = '  +4.41%  '

$ws.Range('E34').Value = '  +2.93%  '

$ws.Range('E35').Value = '  +3.71%  '

$ws.Range('D36').Formula = "'" + '0.957'
$ws.Range('E36').Value = '  +3.09%  '

$ws.Range('E37').Value = '  +0.47%  '

$ws.Range('E38').Value = '  +0.91%  '

$ws.Range('D39').Formula = "'" + '0.522'
$ws.Range('E39').Value = '  +0.45%  '

$ws.Range('E40').Value = '  +1.07%  '

$ws.Range('E41').Value = '  +0.41%  '

$ws.Range('D42').Formula = "'" + '5.57'
$ws.Range('E42').Value = '  -2.73%  '

$ws.Range('D43').Formula = "'" + '0.989'
$ws.Range('E43').Value = '  -0.42%  '

$ws.Range('D44').Formula = "'" + '2.27'
$ws.Range('E44').Value = '  +3.58%  '

$ws.Range('D45').Formula = "'" + '63.57'
$ws.Range('E45').Value = '  +1.98%  '

$ws.Range('D46').Formula = "'" + '1.75'
$ws.Range('E46').Value = '  +0.82%  '

$ws.Range('D47').Value = '1.687.10'
$ws.Range('E47').Value = '  +0.82%  '

$ws.Range('D48').Formula = "'" + '86.18'
$ws.Range('E48').Value = '  +0.48%  '

$ws.Range('D49').Formula = "'" + '0.0517'
$ws.Range('E49').Value = '  +1.44%  '

$ws.Range('D50').Formula = "'" + '0.0956'
$ws.Range('E50').Value = '  +0.98%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Formula = "'" + '1.01'
$ws.Range('E51').Value = '  +0.32%  '
